$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.004605400000000093
$ws.Range("E2").Value = 0.3696322451519838
$ws.Range("I2").Value = 0.388768051676684
$ws.Range("L2").Value = 0.5736454999999999
$ws.Range("M2").Value = 0.082778
$ws.Range("N2").Value = 12.96177663458551
$ws.Range("O2").Value = 3.39354415621619

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.009315151539373345
$ws.Range("B2").Value = 0.07940568642667734
$ws.Range("E2").Value = 0.35180893228384
$ws.Range("I2").Value = 0.7374698220842428
$ws.Range("L2").Value = 0.1242953967812489
$ws.Range("M2").Value = 0.08217333333333333
$ws.Range("N2").Value = 8.857386976801159
$ws.Range("O2").Value = 3.769864170149332

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.1014415319200664
$ws.Range("B2").Value = 0.01832281100914623
$ws.Range("E2").Value = 0.1561358999999998
$ws.Range("I2").Value = 0.4339826273285834
$ws.Range("M2").Value = 0.05117433794284604
$ws.Range("N2").Value = 8.609687601718868
$ws.Range("O2").Value = 5.540089747764068
